$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the first occurrence of $old (search text) with $new inside
# a narrow range right after an anchor phrase. Using wdReplaceOne (1) instead
# of wdReplaceAll (2) keeps the substitution confined to the supplied range.
# ---------------------------------------------------------------------------
function Replace-AfterAnchor($anchorText, $old, $new, $slack) {
    $a = $d.Content
    $a.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $scope = $d.Range($a.End, $a.End + $slack)
    $scope.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 1) | Out-Null
}

# 1 & 3: "English" appears twice (language-switch links); both become the
# same Thai word, so a simple replace-all on the whole document is safe.
$d.Content.Find.Execute("English", $true, $false, $false, $false, $false, $true, 1, $false, "ภาษาอังกฤษ", 2) | Out-Null

# 2: language list next to the hyperlink
$d.Content.Find.Execute(" / Portuguese / French / Thai / Vietnamese / Spanish", $true, $false, $false, $false, $false, $true, 1, $false, " / ภาษาโปรตุเกส / ภาษาฝรั่งเศส /ภาษาไทย / ภาษาเวียดนาม / ภาษาสเปน", 2) | Out-Null

# 4: "Brief" heading
$d.Content.Find.Execute("Brief", $true, $false, $false, $false, $false, $true, 1, $false, "บทย่อ", 2) | Out-Null

# 5: brief description paragraph
$d.Content.Find.Execute("An email sent upon verification to partners in the target country who have sent the correct documents. It will be sent via customer.io", $true, $false, $false, $false, $false, $true, 1, $false, "An email sent upon verification to partners in the target country who have sent the correct documents. โดยมันจะถูกส่งผ่านทาง customer.io", 2) | Out-Null

# 6: "Target audience" heading
$d.Content.Find.Execute("Target audience", $true, $false, $false, $false, $false, $true, 1, $false, "กลุ่มเป้าหมาย", 2) | Out-Null

# 7: subject/heading line
$d.Content.Find.Execute("Your documents have been verified!", $true, $false, $false, $false, $false, $true, 1, $false, "เอกสารของคุณได้รับการตรวจสอบยืนยันแล้ว!", 2) | Out-Null

# 8: greeting
$d.Content.Find.Execute("Hi ", $true, $false, $false, $false, $false, $true, 1, $false, "สวัสดี ", 2) | Out-Null

# 9: ", " right after [PARTNER NAME] -> " "
Replace-AfterAnchor "[PARTNER NAME]" ", " " " 6

# 10: "We've reviewed..." lead-in
$d.Content.Find.Execute("We’ve reviewed the documents you’ve sent us for the ", $true, $false, $false, $false, $false, $true, 1, $false, "เราได้ตรวจสอบเอกสารที่คุณส่งมาให้เราสำหรับงาน ", 2) | Out-Null

# 11: "... and all of them have been verified! "
$d.Content.Find.Execute(" and all of them have been verified! ", $true, $false, $false, $false, $false, $true, 1, $false, " และเอกสารทั้งหมดได้รับการตรวจสอบยืนยันเรียบร้อยแล้ว! ", 2) | Out-Null

# 12: "We'll be sending..." paragraph
$d.Content.Find.Execute("We’ll be sending out more details about the event soon, including the agenda and travel arrangements, so make sure to check your inbox regularly.", $true, $false, $false, $false, $false, $true, 1, $false, "เราจะส่งรายละเอียดเพิ่มเติมเกี่ยวกับกิจกรรมไปให้คุณในเร็วๆ นี้ รวมถึงกำหนดการและการเตรียมเรื่องการเดินทาง ดังนั้นโปรดตรวจดูกล่องข้อความอีเมล์ของคุณอย่างสม่ำเสมอ", 2) | Out-Null

# 13: "If you have any questions, please contact us via "
$d.Content.Find.Execute("If you have any questions, please contact us via ", $true, $false, $false, $false, $false, $true, 1, $false, "หากคุณมีคำถามใดๆ กรุณาติดต่อเราผ่านทาง ", 2) | Out-Null

# 14: "live chat" hyperlink text
$d.Content.Find.Execute("live chat", $true, $false, $false, $false, $false, $true, 1, $false, "แชทสด", 2) | Out-Null

# 15: " or " between "live chat" and "WhatsApp" hyperlinks -> " หรือทาง "
Replace-AfterAnchor "แชทสด" " or " " หรือทาง " 10

# 16: ". " right after the "WhatsApp" hyperlink -> " "
Replace-AfterAnchor "WhatsApp" ". " " " 6

# 17: "If you have any questions, please contact your country manager, "
$d.Content.Find.Execute("If you have any questions, please contact your country manager, ", $true, $false, $false, $false, $false, $true, 1, $false, "หากคุณมีคำถามใดๆ โปรดติดต่อผู้จัดการประจำประเทศของคุณซึ่งได้แก่ ", 2) | Out-Null

# 18: ", at " -> " ที่ "
$d.Content.Find.Execute(", at ", $true, $false, $false, $false, $false, $true, 1, $false, " ที่ ", 2) | Out-Null

# 19: " or " between [EMAIL ADDRESS] and [WHATSAPP NO] -> " หรือ "
Replace-AfterAnchor "[EMAIL ADDRESS]" " or " " หรือ " 10

# 20: " (WhatsApp). " -> " (WhatsApp) "
$d.Content.Find.Execute(" (WhatsApp). ", $true, $false, $false, $false, $false, $true, 1, $false, " (WhatsApp) ", 2) | Out-Null

# 21: the margin comment text
$c = $d.Comments.Item(1)
$c.Range.Text = "เลือกอย่างใดอย่างหนึ่ง"
